$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "followup" -> "post-intervention" wording in the class definition
# for rows 82-87 (column I).
$newText = "A measurement scale used to measure post-intervention data in a study."
$ws.Range("I82").Value = $newText
$ws.Range("I83").Value = $newText
$ws.Range("I84").Value = $newText
$ws.Range("I85").Value = $newText
$ws.Range("I86").Value = $newText
$ws.Range("I87").Value = $newText

# Remove the blank row 669 (a stray empty row in the middle of the
# "LSR1_for_ontology.csv" block), shifting all following rows up by one
# and shrinking the used range from N685 down to N684.
$ws.Rows(669).Delete()
